$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to remove the year suffix
$ws.Name = "g3.11b média"

# Add new "Ano" column (D) with a header matching the style of the other headers
$ws.Range("D1").Value = "Ano"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill the "Ano" column with the period value for all data rows
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 4).Value = "2013-2023"
}
